$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 68.5
$ws.Range("I5").Value = 68.5
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 68.5
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 46.5
$ws.Range("N5").Value = ""
$ws.Range("H12").Value = 416.81818
$ws.Range("J12").Value = 487.14285
$ws.Range("L12").Value = 487.14285
$ws.Range("N12").Value = -827.14285
$ws.Range("H74").Value = 3949.4
$ws.Range("I74").Value = 3860.4443
$ws.Range("K74").Value = 3860.4443
$ws.Range("M74").Value = -2924.4443
$ws.Range("H77").Value = 3949.4
$ws.Range("I77").Value = 3860.4443
$ws.Range("K77").Value = 19302.2215
$ws.Range("M77").Value = -14622.2215
$ws.Range("H82").Value = 20858.2
$ws.Range("I82").Value = 15398
$ws.Range("K82").Value = 46194
$ws.Range("M82").Value = -45788
$ws.Range("H85").Value = 20858.2
$ws.Range("I85").Value = 15398
$ws.Range("K85").Value = 46194
$ws.Range("M85").Value = -44790
$ws.Range("H86").Value = 5331.75
$ws.Range("J86").Value = 5527.6665
$ws.Range("L86").Value = 5527.6665
$ws.Range("N86").Value = -7773.6665
$ws.Range("H89").Value = 5331.75
$ws.Range("J89").Value = 5527.6665
$ws.Range("L89").Value = 27638.3325
$ws.Range("N89").Value = -38870.3325
$ws.Range("H98").Value = 83335210
$ws.Range("I98").Value = 100001650
$ws.Range("J98").Value = 2999.5
$ws.Range("K98").Value = 100001650
$ws.Range("L98").Value = 2999.5
$ws.Range("M98").Value = -100000152
$ws.Range("N98").Value = -5995.5
$ws.Range("H122").Value = 83335210
$ws.Range("I122").Value = 100001650
$ws.Range("J122").Value = 2999.5
$ws.Range("K122").Value = 300004950
$ws.Range("L122").Value = 8998.5
$ws.Range("M122").Value = -300002500
$ws.Range("N122").Value = -13898.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 28682.15
$ws.Range("I132").Value = 29739.445
$ws.Range("K132").Value = 89218.33499999999
$ws.Range("M132").Value = -86688.33499999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6299.6294
$ws.Range("I20").Value = 5939.5
$ws.Range("K20").Value = 5939.5
$ws.Range("M20").Value = -5692.5
$ws.Range("H22").Value = 446.1111
$ws.Range("I22").Value = 409.42856
$ws.Range("J22").Value = 574.5
$ws.Range("K22").Value = 409.42856
$ws.Range("L22").Value = 574.5
$ws.Range("M22").Value = -236.42856
$ws.Range("N22").Value = -920.5
$ws.Range("H94").Value = 1917.1333
$ws.Range("I94").Value = 666.6667
$ws.Range("J94").Value = 2750.7778
$ws.Range("K94").Value = 666.6667
$ws.Range("L94").Value = 2750.7778
$ws.Range("M94").Value = -215.6667
$ws.Range("N94").Value = -3652.7778

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 678417.6
$ws.Range("J31").Value = 1071472.9
$ws.Range("L31").Value = 1071472.9
$ws.Range("N31").Value = -1072062.9
$ws.Range("H34").Value = 678417.6
$ws.Range("J34").Value = 1071472.9
$ws.Range("L34").Value = 1071472.9
$ws.Range("N34").Value = -1071876.9

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 171.46666
$ws.Range("I2").Value = 80.333336
$ws.Range("J2").Value = 181.59259
$ws.Range("K2").Value = 482.000016
$ws.Range("L2").Value = 1089.55554
$ws.Range("M2").Value = -369.000016
$ws.Range("N2").Value = -1315.55554
$ws.Range("H37").Value = 61356.145
$ws.Range("J37").Value = 61356.145
$ws.Range("L37").Value = 184068.435
$ws.Range("N37").Value = -184292.435
$ws.Range("H86").Value = 4380
$ws.Range("J86").Value = 5250
$ws.Range("L86").Value = 15750
$ws.Range("N86").Value = -18122
$ws.Range("H89").Value = 4380
$ws.Range("J89").Value = 5250
$ws.Range("L89").Value = 47250
$ws.Range("N89").Value = -59106
$ws.Range("H94").Value = 3999.923
$ws.Range("H131").Value = 17593.428
$ws.Range("J131").Value = 20410
$ws.Range("L131").Value = 61230
$ws.Range("N131").Value = -71310
$ws.Range("H133").Value = 3132.5
$ws.Range("I133").Value = 3132.5
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 9397.5
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -4337.5
$ws.Range("N133").Value = ""

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5064.1665
$ws.Range("I80").Value = 3296.3333
$ws.Range("J80").Value = 6832
$ws.Range("K80").Value = 3296.3333
$ws.Range("L80").Value = 6832
$ws.Range("M80").Value = -2298.3333
$ws.Range("N80").Value = -8828
$ws.Range("H83").Value = 5064.1665
$ws.Range("I83").Value = 3296.3333
$ws.Range("J83").Value = 6832
$ws.Range("K83").Value = 16481.6665
$ws.Range("L83").Value = 34160
$ws.Range("M83").Value = -11489.6665
$ws.Range("N83").Value = -44144
$ws.Range("H136").Value = 11066
$ws.Range("J136").Value = 11066
$ws.Range("L136").Value = 33198
$ws.Range("N136").Value = -38298

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 13886.35
$ws.Range("J7").Value = 14870.583
$ws.Range("L7").Value = 14870.583
$ws.Range("N7").Value = -15094.583
$ws.Range("H99").Value = 38663.332
$ws.Range("I99").Value = 34000
$ws.Range("K99").Value = 34000
$ws.Range("M99").Value = -31005
$ws.Range("H126").Value = 13886.35
$ws.Range("J126").Value = 14870.583
$ws.Range("L126").Value = 44611.749
$ws.Range("N126").Value = -49551.749
$ws.Range("H132").Value = 318676.6
$ws.Range("I132").Value = 6643
$ws.Range("K132").Value = 19929
$ws.Range("M132").Value = -17399

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 25576.54
$ws.Range("I54").Value = 18545.455
$ws.Range("K54").Value = 18545.455
$ws.Range("M54").Value = -18025.455
$ws.Range("H81").Value = 7935.222
$ws.Range("I81").Value = 3609.6
$ws.Range("J81").Value = 13342.25
$ws.Range("K81").Value = 7219.2
$ws.Range("L81").Value = 26684.5
$ws.Range("M81").Value = -6158.2
$ws.Range("N81").Value = -28806.5
$ws.Range("H84").Value = 7935.222
$ws.Range("I84").Value = 3609.6
$ws.Range("J84").Value = 13342.25
$ws.Range("K84").Value = 36096
$ws.Range("L84").Value = 133422.5
$ws.Range("M84").Value = -30792
$ws.Range("N84").Value = -144030.5
